$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("CONCEPTO_TECNICO_NO_OPER", $true, $false, $false, $false, $false, $true, 1, $false, "CONCEPTO_TECNICO_NO_OPER", 2)
